$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume table update (GitHub Actions data refresh)
# Column D ("Price") values that are purely numeric text (e.g. "310.90")
# must be force-formatted as Text first, otherwise Excel auto-converts
# them to numbers (stripping the trailing zero) when assigned via .Value.

$ws.Range("D2").Value = '43.446.72'
$ws.Range("E2").Value = '  +2.82%  '
$ws.Range("D3").Value = '2.309.84'
$ws.Range("E3").Value = '  +1.80%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.90'
$ws.Range("E5").Value = '  +1.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '101.86'
$ws.Range("E6").Value = '  +4.72%  '
$ws.Range("E7").Value = '  +1.58%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  +7.68%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.73'
$ws.Range("E10").Value = '  +1.64%  '
$ws.Range("E11").Value = '  +3.09%  '
$ws.Range("E12").Value = '  -0.55%  '
$ws.Range("E13").Value = '  +0.43%  '
$ws.Range("D14").Value = '2.667.74'
$ws.Range("E14").Value = '  +1.74%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.95'
$ws.Range("E15").Value = '  +1.39%  '
$ws.Range("D16").Value = '2.311.89'
$ws.Range("E16").Value = '  +1.62%  '
$ws.Range("E17").Value = '  +2.09%  '
$ws.Range("D18").Value = '43.349.44'
$ws.Range("E18").Value = '  +2.95%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.34'
$ws.Range("E19").Value = '  -0.03%  '
$ws.Range("E20").Value = '  +2.42%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.16'
$ws.Range("E21").Value = '  +2.38%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.08'
$ws.Range("E22").Value = '  +0.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '241.32'
$ws.Range("E23").Value = '  +1.42%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.02'
$ws.Range("E24").Value = '  +3.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.61'
$ws.Range("E25").Value = '  +1.47%  '
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("E27").Value = '  -1.54%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '24.62'
$ws.Range("E28").Value = '  +4.42%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '36.84'
$ws.Range("E29").Value = '  -2.92%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.65'
$ws.Range("E30").Value = '  +1.26%  '
$ws.Range("E31").Value = '  +0.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '167.22'
$ws.Range("E32").Value = '  +3.81%  '
$ws.Range("E33").Value = '  +0.87%  '
$ws.Range("E34").Value = '  -0.02%  '
$ws.Range("E35").Value = '  +0.94%  '
$ws.Range("E36").Value = '  +5.60%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.06'
$ws.Range("E37").Value = '  -2.56%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '17.58'
$ws.Range("E38").Value = '  -0.44%  '
$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.88'
$ws.Range("E39").Value = '  +3.08%  '
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.106'
$ws.Range("E40").Value = '  +1.24%  '
$ws.Range("E41").Value = '  +1.39%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.35'
$ws.Range("E42").Value = '  +7.06%  '
$ws.Range("E43").Value = '  -1.03%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.33'
$ws.Range("E44").Value = '  +0.07%  '
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0289'
$ws.Range("E45").Value = '  +2.78%  '
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '1.967.72'
$ws.Range("E46").Value = '  +1.32%  '
$ws.Range("E47").Value = '  +2.41%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.96'
$ws.Range("E48").Value = '  +0.73%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '55.54'
$ws.Range("E49").Value = '  +3.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.91'
$ws.Range("E50").Value = '  +5.39%  '
$ws.Range("E51").Value = '  +6.99%  '
